$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.015.33"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "1.987.10"
$ws.Range("E3").Value = "  -2.48%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -6.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.602"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.61%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.86"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -5.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.372"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.70"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.41%  "
$ws.Range("E11").Value = "  -5.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0976"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.72%  "
$ws.Range("D13").Value = "2.277.61"
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.03"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.75"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.755"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -8.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.03"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -6.65%  "
$ws.Range("D18").Value = "1.990.51"
$ws.Range("E18").Value = "  -2.60%  "
$ws.Range("D19").Value = "36.957.03"
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.99"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.11%  "
$ws.Range("D21").Value = "0.0₃0808"
$ws.Range("E21").Value = "  -5.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.41"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("E23").Value = "  -5.10%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -9.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.34"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.80"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.62"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.08"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.87%  "
$ws.Range("E30").Value = "  -11.35%  "
$ws.Range("E31").Value = "  -3.40%  "
$ws.Range("E32").Value = "  -3.67%  "
$ws.Range("E33").Value = "  -7.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0609"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -9.08%  "
$ws.Range("E35").Value = "  -7.54%  "
$ws.Range("E36").Value = "  -6.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("E39").Value = "  -4.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.21"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.73%  "
$ws.Range("E41").Value = "  +2.36%  "
$ws.Range("D42").Value = "1.430.46"
$ws.Range("E42").Value = "  +1.71%  "
$ws.Range("E43").Value = "  -6.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.11"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0875"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -9.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.04"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.11"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -7.16%  "
$ws.Range("E48").Value = "  -5.16%  "
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.63"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +12.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.64"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -10.28%  "
